$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ML Project Checklist")

# Fill in "Y" for the Done column (C) on rows 6-10, matching the
# already-populated cells in rows 2-5.
$ws.Range("C6").Value = "Y"
$ws.Range("C7").Value = "Y"
$ws.Range("C8").Value = "Y"
$ws.Range("C9").Value = "Y"
$ws.Range("C10").Value = "Y"

# Update the active selection to C10 (was D5).
$ws.Range("C10").Select()
